$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values for rows that changed ---

# Row 10
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '8870322 - Fabiano Fernandes Bargos'
$ws.Range('C10').Value = '8870322 - Fabiano Fernandes Bargos'

# Row 13
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'

# Row 14
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = 'Introduction to Matlab (or SciLab, Octave, Freemat, etc.), roots of nonlinear equations, systems of equations, least-squares fitting of curves to data, numerical integration, and solving ordinary differential equations.'
$ws.Range('C14').Value = 'Introduction to Matlab (or SciLab, Octave, Freemat, etc.), roots of nonlinear equations, systems of equations, least-squares fitting of curves to data, numerical integration, and solving ordinary differential equations.'

# Row 15
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '01/01/2018'
$ws.Range('C15').Value = '01/01/2018'

# Row 16
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B16').Value = '1.Introduction to Matlab (or SciLab, Octave, Freemat, etc.); Simple calculations; Use of variables and functions. Matrices and Vectors; Plotting.2.Matlab Programming; Script m-files, Function m-files, input and output, flow control, vectorization, global variables.3.Finding Roots of Nonlinear Equations: Fixed point, bisection, and Newton’s method.4.Linear Algebra (matrix and vector properties and operations ).5.Solving systems of linear equations.6.Least-squares fitting of curves to data.7.Interpolation.8.Numerical Integration.9.Solving Ordinary Differential Equations'
$ws.Range('C16').Value = '1.Introduction to Matlab (or SciLab, Octave, Freemat, etc.); Simple calculations; Use of variables and functions. Matrices and Vectors; Plotting.2.Matlab Programming; Script m-files, Function m-files, input and output, flow control, vectorization, global variables.3.Finding Roots of Nonlinear Equations: Fixed point, bisection, and Newton’s method.4.Linear Algebra (matrix and vector properties and operations ).5.Solving systems of linear equations.6.Least-squares fitting of curves to data.7.Interpolation.8.Numerical Integration.9.Solving Ordinary Differential Equations'

# Row 17
$ws.Range('A17').Value = 'Avaliação:'

# Row 18
$ws.Range('A18').Value = 'Método:'
$ws.Range('B18').Value = '8870322 - Fabiano Fernandes Bargos'
$ws.Range('C18').Value = '8870322 - Fabiano Fernandes Bargos'

# Row 19
$ws.Range('A19').Value = 'Critério:'
$ws.Range('B19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range('C19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# Row 20
$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B20').Value = 'NF≥ 5,0.'
$ws.Range('C20').Value = 'NF≥ 5,0.'

# Row 21
$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B21').Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range('C21').Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'

# Row 22
$ws.Range('A22').Value = 'Requisitos:'

# Row 23
$ws.Range('B23').Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'
$ws.Range('C23').Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'

# Row 24
$ws.Range('B24').Value = 'LOB1036 -  Geometria Analítica  (Requisito fraco)
'
$ws.Range('C24').Value = 'LOB1036 -  Geometria Analítica  (Requisito fraco)
'

# --- Clear B/C cells on rows that no longer have that content (row shrinks to label-only) ---
$ws.Range('B17:C17').Clear()
$ws.Range('B22:C22').Clear()

# --- Adjust row heights to match final layout ---
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).EntireRow.AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).EntireRow.AutoFit()
$ws.Rows(23).RowHeight = 30

# --- Remove now-obsolete last row (content has shifted up into row 24) ---
$ws.Rows(25).Delete()
